# Cập nhật API Login
# - Split the old "login fail" JSON sample on the Login sheet into two
#   distinct error responses: "wrong email" and "wrong password".
# - Make the Login sheet the active tab/selection (was List).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Login")

# Insert a new row 24 (inherits formatting/style from row 23 above it)
# so B21:B24 become the new "{ / error_code / message / }" JSON block.
$ws.Rows(24).Insert()

# New second JSON error example (wrong password) placed in the
# previously-empty B21:B24 block.
$ws.Range("B21").Value = "{"
$ws.Range("B22").Value = " ""error_code"" : ""1"","
$ws.Range("B23").Value = " ""message"" : ""wrong password"""
$ws.Range("B24").Value = "}"
$ws.Rows(24).RowHeight = 15.6

# The first JSON error example's message changes from the generic
# "login fail" to the more specific "wrong email".
$ws.Range("B18").Value = " ""message"" : ""wrong email"""

# Login becomes the active sheet/tab (List was previously active),
# with the on-screen selection moved to H12.
$ws.Activate()
$ws.Range("H12").Select()
